$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameter")

# Row 4 of the "parameter" sheet describes the "plan.loglevel" parameter;
# repurpose it into the new "allowsplits" parameter (default: false).
$ws.Range("A4").Value = "allowsplits"

# Column B holds the default value. We need the literal text "true" (not a
# Boolean) because this column is a shared-string "value" column, so write
# it through a temp formula + PasteSpecial(values) round-trip, which keeps
# the literal text instead of auto-coercing "true" to a Boolean.
$ws.Range("Z1").Formula = "=""true"""
$ws.Range("Z1").Copy()
$ws.Range("B4").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("C4").Value = "Controls whether a sales order or forecast can be split across multiple manufacturing orders during planning. Default: false"

# Move the active sheet/selection onto the parameter sheet, row 4 selected.
$ws.Activate()
$ws.Range("A4:XFD4").Select()
